$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume }
# (only the cells that actually changed per the diff are listed per row)
$updates = @{
    2  = @{ D = "261.07";     E = "1.64%" }
    3  = @{ D = "27.35";      E = "1.74%" }
    4  = @{ D = "4.707";      E = "1.48%" }
    5  = @{ D = "0.06078";    E = "3.14%" }
    6  = @{ D = "6.674";      E = "0.99%" }
    7  = @{ E = "-0.48%" }
    8  = @{ D = "0.9243";     E = "0.42%" }
    9  = @{ D = "0.1402";     E = "1.64%" }
    10 = @{ D = "0.04953";    E = "15.55%" }
    11 = @{ D = "0.07136";    E = "1.98%" }
    12 = @{ D = "0.03098";    E = "1.32%" }
    13 = @{ D = "0.09068";    E = "-0.39%" }
    14 = @{ D = "0.001537";   E = "0.71%" }
    15 = @{ D = "0.0006094";  E = "0.71%" }
    16 = @{ D = "0.006108";   E = "1.57%" }
    17 = @{ E = "-0.54%" }
    18 = @{ D = "3.143";      E = "-0.63%" }
    19 = @{ E = "-0.33%" }
    20 = @{ E = "2.34%" }
    21 = @{ D = "0.1288";     E = "-0.49%" }
    22 = @{ D = "4.090";      E = "4.75%" }
    23 = @{ D = "0.04233";    E = "-0.58%" }
    24 = @{ D = "0.001222";   E = "-0.07%" }
    25 = @{ E = "-8.85%" }
    26 = @{ E = "0.00%" }
    27 = @{ D = "0.0001575";  E = "3.40%" }
    40 = @{ D = "0.03867";    E = "2.33%" }
    41 = @{ E = "1.22%" }
    42 = @{ D = "0.004093";   E = "5.45%" }
    43 = @{ D = "0.01636";    E = "21.33%" }
    44 = @{ D = "0.002217";   E = "-8.04%" }
    45 = @{ D = "0.00005140"; E = "-4.12%" }
    46 = @{ E = "0.04%" }
    47 = @{ D = "0.05446";    E = "19.70%" }
    48 = @{ E = "-46.44%" }
    49 = @{ E = "0.04%" }
    50 = @{ E = "0.04%" }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        # These columns store plain text (e.g. "261.07", "1.64%") in the
        # source workbook, not numbers. Force text interpretation so Excel
        # doesn't silently convert the assigned string into a numeric value,
        # then restore the default "Normal" style so no stray number-format
        # style gets attached to the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col]
        $cell.Style = "Normal"
    }
}
